# Split the combined "address: ico" and "address phone" lines in the
# contract header into two separate lines, each terminated with a
# manual line break (w:br), matching the target diff.

$d = $word.ActiveDocument

# 1) "Sídlo: Nad Kampusem 821/4, [[ADDRESS_1]]: [[ICO_1]]"
#    -> "Sídlo: Nad Kampusem 821/4, [[ADDRESS_1]]" + line break + "[[ICO_1]]"
$d.Content.Find.Execute(
    "Sídlo: Nad Kampusem 821/4, [[ADDRESS_1]]: [[ICO_1]]",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Sídlo: Nad Kampusem 821/4, [[ADDRESS_1]]^l[[ICO_1]]",
    2
)

# 2) "[[ADDRESS_2]][[PHONE_2]]"
#    -> "[[ADDRESS_2]]" + line break + "[[PHONE_2]]"
$d.Content.Find.Execute(
    "[[ADDRESS_2]][[PHONE_2]]",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "[[ADDRESS_2]]^l[[PHONE_2]]",
    2
)
